$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 25 (2025-12) stats
$ws.Range("B25").Value = 6440
$ws.Range("C25").Value = 1004
$ws.Range("D25").Value = 5999261
$ws.Range("E25").Value = 931.5622670807453
$ws.Range("F25").Value = 9.319300628076732
$ws.Range("G25").Value = 7.036247334754808
$ws.Range("H25").Value = 25.64302597403141
